$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on sheet "Hoja1" (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Value2
$text = $text -replace [regex]::Escape("✅ 1000 Bs = 7.5 = 30400.28 pesos"), "✅ 1000 Bs = 7.43 = 30171.0 pesos"
$text = $text -replace [regex]::Escape("✅ 30400.28 pesos = 7.48 = 968.01 Bs"), "✅ 30171.0 pesos = 7.41 = 956.28 Bs"
$cellA1.Value2 = $text

# --- Update the tasas table on sheet "tasas" ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value2 = 134.5
$wsTasas.Range("O10").Value2 = 4058
$wsTasas.Range("N12").Value2 = 4070
$wsTasas.Range("O12").Value2 = 129
